$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the table with the 2024/11/28 column (new column CC = index 81) ---

# Give column CC the same width as the rest of the data columns (stored OOXML width 12).
# Excels ColumnWidth property reports a value that is offset from the raw stored width,
# so a ColumnWidth of (11 + 1/6) serializes to width="12" in the sheet XML.
$ws.Columns.Item(81).ColumnWidth = 11 + 1/6

# Reference cells that already carry each of the value-based styles used throughout the sheet,
# used as PasteSpecial(xlPasteFormats) sources so that the new cells reuse the exact same
# style records instead of creating new (duplicate) ones:
#  B1 -> s=1 (header row style, font "Meiryo", General format)
#  B2 -> s=1 (normal,            value >= 140)
#  N2 -> s=3 (light-blue fill,    125 <= value < 140)
#  D2 -> s=2 (yellow fill,        value < 125)

# --- Row 1 header: the date label "2024/11/28" (stored as text, like the other date headers) ---
$ws.Range("CC1").NumberFormat = "@"
$ws.Range("CC1").Value = "2024/11/28"
$ws.Range("B1").Copy()
$ws.Range("CC1").PasteSpecial(-4122)

# --- Data rows 2-53 ---
$ws.Range("B2").Copy()
$ws.Range("CC2").PasteSpecial(-4122)
$ws.Range("CC2").Value = 198.9

$ws.Range("N2").Copy()
$ws.Range("CC3").PasteSpecial(-4122)
$ws.Range("CC3").Value = 137.2

$ws.Range("B2").Copy()
$ws.Range("CC4").PasteSpecial(-4122)
$ws.Range("CC4").Value = 231.8

$ws.Range("B2").Copy()
$ws.Range("CC5").PasteSpecial(-4122)
$ws.Range("CC5").Value = 184.3

$ws.Range("B2").Copy()
$ws.Range("CC6").PasteSpecial(-4122)
$ws.Range("CC6").Value = 202.1

$ws.Range("B2").Copy()
$ws.Range("CC7").PasteSpecial(-4122)
$ws.Range("CC7").Value = 209.6

$ws.Range("B2").Copy()
$ws.Range("CC8").PasteSpecial(-4122)
$ws.Range("CC8").Value = 140.1

$ws.Range("B2").Copy()
$ws.Range("CC9").PasteSpecial(-4122)
$ws.Range("CC9").Value = 174.4

$ws.Range("B2").Copy()
$ws.Range("CC10").PasteSpecial(-4122)
$ws.Range("CC10").Value = 184.1

$ws.Range("B2").Copy()
$ws.Range("CC11").PasteSpecial(-4122)
$ws.Range("CC11").Value = 142.3

$ws.Range("B2").Copy()
$ws.Range("CC12").PasteSpecial(-4122)
$ws.Range("CC12").Value = 203.7

$ws.Range("B2").Copy()
$ws.Range("CC13").PasteSpecial(-4122)
$ws.Range("CC13").Value = 151.7

$ws.Range("B2").Copy()
$ws.Range("CC14").PasteSpecial(-4122)
$ws.Range("CC14").Value = 140.5

$ws.Range("B2").Copy()
$ws.Range("CC15").PasteSpecial(-4122)
$ws.Range("CC15").Value = 225.8

$ws.Range("N2").Copy()
$ws.Range("CC16").PasteSpecial(-4122)
$ws.Range("CC16").Value = 138.3

$ws.Range("B2").Copy()
$ws.Range("CC17").PasteSpecial(-4122)
$ws.Range("CC17").Value = 190.3

$ws.Range("B2").Copy()
$ws.Range("CC18").PasteSpecial(-4122)
$ws.Range("CC18").Value = 182.6

$ws.Range("N2").Copy()
$ws.Range("CC19").PasteSpecial(-4122)
$ws.Range("CC19").Value = 138.2

$ws.Range("B2").Copy()
$ws.Range("CC20").PasteSpecial(-4122)
$ws.Range("CC20").Value = 162.5

$ws.Range("B2").Copy()
$ws.Range("CC21").PasteSpecial(-4122)
$ws.Range("CC21").Value = 168.7

$ws.Range("N2").Copy()
$ws.Range("CC22").PasteSpecial(-4122)
$ws.Range("CC22").Value = 127.2

$ws.Range("N2").Copy()
$ws.Range("CC23").PasteSpecial(-4122)
$ws.Range("CC23").Value = 128.6

$ws.Range("N2").Copy()
$ws.Range("CC24").PasteSpecial(-4122)
$ws.Range("CC24").Value = 129.9

$ws.Range("B2").Copy()
$ws.Range("CC25").PasteSpecial(-4122)
$ws.Range("CC25").Value = 174.8

$ws.Range("B2").Copy()
$ws.Range("CC26").PasteSpecial(-4122)
$ws.Range("CC26").Value = 165.7

$ws.Range("B2").Copy()
$ws.Range("CC27").PasteSpecial(-4122)
$ws.Range("CC27").Value = 186.2

$ws.Range("B2").Copy()
$ws.Range("CC28").PasteSpecial(-4122)
$ws.Range("CC28").Value = 141.7

$ws.Range("B2").Copy()
$ws.Range("CC29").PasteSpecial(-4122)
$ws.Range("CC29").Value = 161.6

$ws.Range("N2").Copy()
$ws.Range("CC30").PasteSpecial(-4122)
$ws.Range("CC30").Value = 128.8

$ws.Range("B2").Copy()
$ws.Range("CC31").PasteSpecial(-4122)
$ws.Range("CC31").Value = 188

$ws.Range("B2").Copy()
$ws.Range("CC32").PasteSpecial(-4122)
$ws.Range("CC32").Value = 158.9

$ws.Range("B2").Copy()
$ws.Range("CC33").PasteSpecial(-4122)
$ws.Range("CC33").Value = 148.1

$ws.Range("B2").Copy()
$ws.Range("CC34").PasteSpecial(-4122)
$ws.Range("CC34").Value = 365.3

$ws.Range("B2").Copy()
$ws.Range("CC35").PasteSpecial(-4122)
$ws.Range("CC35").Value = 140

$ws.Range("B2").Copy()
$ws.Range("CC36").PasteSpecial(-4122)
$ws.Range("CC36").Value = 155.2

$ws.Range("B2").Copy()
$ws.Range("CC37").PasteSpecial(-4122)
$ws.Range("CC37").Value = 264.5

$ws.Range("N2").Copy()
$ws.Range("CC38").PasteSpecial(-4122)
$ws.Range("CC38").Value = 137.8

$ws.Range("B2").Copy()
$ws.Range("CC39").PasteSpecial(-4122)
$ws.Range("CC39").Value = 239.7

$ws.Range("B2").Copy()
$ws.Range("CC40").PasteSpecial(-4122)
$ws.Range("CC40").Value = 157.7

$ws.Range("B2").Copy()
$ws.Range("CC41").PasteSpecial(-4122)
$ws.Range("CC41").Value = 151.2

$ws.Range("B2").Copy()
$ws.Range("CC42").PasteSpecial(-4122)
$ws.Range("CC42").Value = 187.7

$ws.Range("B2").Copy()
$ws.Range("CC43").PasteSpecial(-4122)
$ws.Range("CC43").Value = 154.8

$ws.Range("B2").Copy()
$ws.Range("CC44").PasteSpecial(-4122)
$ws.Range("CC44").Value = 162.8

$ws.Range("N2").Copy()
$ws.Range("CC45").PasteSpecial(-4122)
$ws.Range("CC45").Value = 134.7

$ws.Range("B2").Copy()
$ws.Range("CC46").PasteSpecial(-4122)
$ws.Range("CC46").Value = 165.4

$ws.Range("B2").Copy()
$ws.Range("CC47").PasteSpecial(-4122)
$ws.Range("CC47").Value = 154.3

$ws.Range("N2").Copy()
$ws.Range("CC48").PasteSpecial(-4122)
$ws.Range("CC48").Value = 136.3

$ws.Range("B2").Copy()
$ws.Range("CC49").PasteSpecial(-4122)
$ws.Range("CC49").Value = 177.9

$ws.Range("B2").Copy()
$ws.Range("CC50").PasteSpecial(-4122)
$ws.Range("CC50").Value = 148.3

$ws.Range("B2").Copy()
$ws.Range("CC51").PasteSpecial(-4122)
$ws.Range("CC51").Value = 194.9

$ws.Range("B2").Copy()
$ws.Range("CC52").PasteSpecial(-4122)
$ws.Range("CC52").Value = 166.8

$ws.Range("B2").Copy()
$ws.Range("CC53").PasteSpecial(-4122)
$ws.Range("CC53").Value = 152.7

